# Auto-generated Excel COM-interop script
# Applies scheduled-runner market price refresh to Asura_Profits data across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2733.1667
$ws.Range("I62").Value = 2150
$ws.Range("J62").Value = 3024.75
$ws.Range("K62").Value = 2150
$ws.Range("L62").Value = 3024.75
$ws.Range("M62").Value = -1526
$ws.Range("N62").Value = -4272.75

$ws.Range("H64").Value = 3312.4849
$ws.Range("I64").Value = 3159.7856
$ws.Range("K64").Value = 3159.7856
$ws.Range("M64").Value = -2911.7856

$ws.Range("H65").Value = 2733.1667
$ws.Range("I65").Value = 2150
$ws.Range("J65").Value = 3024.75
$ws.Range("K65").Value = 10750
$ws.Range("L65").Value = 15123.75
$ws.Range("M65").Value = -7630
$ws.Range("N65").Value = -21363.75

$ws.Range("H67").Value = 3312.4849
$ws.Range("I67").Value = 3159.7856
$ws.Range("K67").Value = 3159.7856
$ws.Range("M67").Value = -2301.7856

$ws.Range("H69").Value = 2500
$ws.Range("I69").Value = 2000
$ws.Range("K69").Value = 6000
$ws.Range("M69").Value = -5126

$ws.Range("H72").Value = 2500
$ws.Range("I72").Value = 2000
$ws.Range("K72").Value = 18000
$ws.Range("M72").Value = -13632

$ws.Range("H80").Value = 6275.263
$ws.Range("I80").Value = 488.85715
$ws.Range("J80").Value = 9650.666999999999
$ws.Range("K80").Value = 1466.57145
$ws.Range("L80").Value = 28952.001
$ws.Range("M80").Value = -468.5714499999999
$ws.Range("N80").Value = -30948.001

$ws.Range("H83").Value = 6275.263
$ws.Range("I83").Value = 488.85715
$ws.Range("J83").Value = 9650.666999999999
$ws.Range("K83").Value = 4399.71435
$ws.Range("L83").Value = 86856.003
$ws.Range("M83").Value = 592.2856499999998
$ws.Range("N83").Value = -96840.003

$ws.Range("H107").Value = 422.33334
$ws.Range("I107").Value = 446.10526
$ws.Range("J107").Value = 332
$ws.Range("K107").Value = 446.10526
$ws.Range("L107").Value = 332
$ws.Range("M107").Value = 1473.89474
$ws.Range("N107").Value = -4172

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 728.61365
$ws.Range("I2").Value = 422.45456
$ws.Range("J2").Value = 1647.091
$ws.Range("K2").Value = 422.45456
$ws.Range("L2").Value = 1647.091
$ws.Range("M2").Value = -309.45456
$ws.Range("N2").Value = -1873.091

$ws.Range("H32").Value = 18530.064
$ws.Range("I32").Value = 19082.512
$ws.Range("K32").Value = 19082.512
$ws.Range("M32").Value = -18795.512

$ws.Range("H45").Value = 2442.2
$ws.Range("I45").Value = 2442.2
$ws.Range("K45").Value = 2442.2
$ws.Range("M45").Value = -2065.2

$ws.Range("H63").Value = 3485.4211
$ws.Range("I63").Value = 3123.5
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 3123.5
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = -2437.5
$ws.Range("N63").Value = -11372

$ws.Range("H66").Value = 3485.4211
$ws.Range("I66").Value = 3123.5
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 15617.5
$ws.Range("L66").Value = 50000
$ws.Range("M66").Value = -12185.5
$ws.Range("N66").Value = -56864

$ws.Range("H88").Value = 2163.875
$ws.Range("I88").Value = 1650
$ws.Range("K88").Value = 1650
$ws.Range("M88").Value = -1244

$ws.Range("H91").Value = 2163.875
$ws.Range("I91").Value = 1650
$ws.Range("K91").Value = 1650
$ws.Range("M91").Value = -246

$ws.Range("H110").Value = 1334.4
$ws.Range("I110").Value = 1447.5454
$ws.Range("J110").Value = 1023.25
$ws.Range("K110").Value = 1447.5454
$ws.Range("L110").Value = 1023.25
$ws.Range("M110").Value = 597.4546
$ws.Range("N110").Value = -5113.25

$ws.Range("H116").Value = 728.61365
$ws.Range("I116").Value = 422.45456
$ws.Range("J116").Value = 1647.091
$ws.Range("K116").Value = 422.45456
$ws.Range("L116").Value = 1647.091
$ws.Range("M116").Value = 1871.54544
$ws.Range("N116").Value = -6235.091

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 728.61365
$ws.Range("I3").Value = 422.45456
$ws.Range("J3").Value = 1647.091
$ws.Range("K3").Value = 422.45456
$ws.Range("L3").Value = 1647.091
$ws.Range("M3").Value = -308.45456
$ws.Range("N3").Value = -1875.091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25643884
$ws.Range("J31").Value = 4955.5
$ws.Range("L31").Value = 4955.5
$ws.Range("N31").Value = -5545.5

$ws.Range("H34").Value = 25643884
$ws.Range("J34").Value = 4955.5
$ws.Range("L34").Value = 4955.5
$ws.Range("N34").Value = -5359.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 6155.8
$ws.Range("J39").Value = 6155.8
$ws.Range("L39").Value = 18467.4
$ws.Range("N39").Value = -19055.4

$ws.Range("H69").Value = 1034.2142
$ws.Range("I69").Value = 500
$ws.Range("J69").Value = 1075.3077
$ws.Range("K69").Value = 1500
$ws.Range("L69").Value = 3225.9231
$ws.Range("M69").Value = -689
$ws.Range("N69").Value = -4847.9231

$ws.Range("H72").Value = 1034.2142
$ws.Range("I72").Value = 500
$ws.Range("J72").Value = 1075.3077
$ws.Range("K72").Value = 4500
$ws.Range("L72").Value = 9677.7693
$ws.Range("M72").Value = -444
$ws.Range("N72").Value = -17789.7693

$ws.Range("H122").Value = 702.1818
$ws.Range("I122").Value = 443.30435
$ws.Range("J122").Value = 1297.6
$ws.Range("K122").Value = 3989.73915
$ws.Range("L122").Value = 11678.4
$ws.Range("M122").Value = -1539.73915
$ws.Range("N122").Value = -16578.4

$ws.Range("H131").Value = 863.59
$ws.Range("I131").Value = 488
$ws.Range("J131").Value = 883.3579
$ws.Range("K131").Value = 1464
$ws.Range("L131").Value = 2650.0737
$ws.Range("M131").Value = 3576
$ws.Range("N131").Value = -12730.0737

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4381
$ws.Range("I80").Value = 5335
$ws.Range("J80").Value = 2950
$ws.Range("K80").Value = 5335
$ws.Range("L80").Value = 2950
$ws.Range("M80").Value = -4337
$ws.Range("N80").Value = -4946

$ws.Range("H83").Value = 4381
$ws.Range("I83").Value = 5335
$ws.Range("J83").Value = 2950
$ws.Range("K83").Value = 26675
$ws.Range("L83").Value = 14750
$ws.Range("M83").Value = -21683
$ws.Range("N83").Value = -24734

$ws.Range("H93").Value = 20250
$ws.Range("J93").Value = 20250
$ws.Range("L93").Value = 20250
$ws.Range("N93").Value = -23994

$ws.Range("H113").Value = 1275.7273
$ws.Range("J113").Value = 1631.25
$ws.Range("L113").Value = 1631.25
$ws.Range("N113").Value = -5971.25

$ws.Range("H123").Value = 10394.772
$ws.Range("J123").Value = 10394.772
$ws.Range("L123").Value = 10394.772
$ws.Range("N123").Value = -15294.772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1522.8572
$ws.Range("I46").Value = 1165
$ws.Range("K46").Value = 1165
$ws.Range("M46").Value = -977

$ws.Range("H61").Value = 14653.177
$ws.Range("I61").Value = 18246.54
$ws.Range("K61").Value = 18246.54
$ws.Range("M61").Value = -18044.54

$ws.Range("H82").Value = 3205.889
$ws.Range("I82").Value = 1990
$ws.Range("J82").Value = 4725.75
$ws.Range("K82").Value = 1990
$ws.Range("L82").Value = 4725.75
$ws.Range("M82").Value = -1629
$ws.Range("N82").Value = -5447.75

$ws.Range("H85").Value = 3205.889
$ws.Range("I85").Value = 1990
$ws.Range("J85").Value = 4725.75
$ws.Range("K85").Value = 1990
$ws.Range("L85").Value = 4725.75
$ws.Range("M85").Value = -742
$ws.Range("N85").Value = -7221.75

$ws.Range("H113").Value = 14653.177
$ws.Range("I113").Value = 18246.54
$ws.Range("K113").Value = 18246.54
$ws.Range("M113").Value = -16076.54

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 80118.38
$ws.Range("J81").Value = 2880
$ws.Range("L81").Value = 5760
$ws.Range("N81").Value = -7882

$ws.Range("H84").Value = 80118.38
$ws.Range("J84").Value = 2880
$ws.Range("L84").Value = 28800
$ws.Range("N84").Value = -39408

$ws.Range("H107").Value = 387.96875
$ws.Range("I107").Value = 268.6087
$ws.Range("J107").Value = 693
$ws.Range("K107").Value = 805.8261
$ws.Range("L107").Value = 2079
$ws.Range("M107").Value = 1114.1739
$ws.Range("N107").Value = -5919

$ws.Range("H113").Value = 688.3871
$ws.Range("I113").Value = 433.42856
$ws.Range("K113").Value = 1300.28568
$ws.Range("M113").Value = 869.71432
